$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 "ui - vuejs" task: status moved from "for now ready" note to "Completed"
$ws.Range("D8").Value = "Приключен"

# Row 9: speedmeter task becomes the speedmeter screen-streaming task
$ws.Range("A9").Value = "speedmeter screen streeaming"
$ws.Range("B9").Value = "screen streaming of the speedmeter app"
$ws.Range("C9").Value = "Радослав, Данило"
$ws.Range("D9").Value = "Приключен"

# Row 10: документация task keeps its text, status becomes "Приключен"
$ws.Range("D10").Value = "Приключен"

# Row 11: new row for the screen-streaming documentation task
$ws.Range("A11").Value = "screen streaming документация"
$ws.Range("B11").Value = "документация относно screen streaming функционлаността"
$ws.Range("B11").WrapText = $true
$ws.Range("C11").Value = "Радослав, Данило"
$ws.Range("D11").Value = "Приключен"

# Drop the leftover stray note that used to live in row 13
$ws.Range("A13").ClearContents()

# Move the active selection to match the edited row
$ws.Range("D11").Select() | Out-Null
